{"js": "// Change 1: header cell \"Result up/Down\" -> \"Result\"\nconst headerResults = context.document.body.search(\"Result up/Down\", { matchCase: true });\nheaderResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < headerResults.items.length; i++) {\n  headerResults.items[i].insertText(\"Result\", Word.InsertLocation.replace);\n}\n\n// Change 2-4: collapse the split \"{{client_visits_X}}\" placeholders (each was\n// spread across three runs: \"{{\", \"client_visits_\", \"X}}\") back into a single\n// run containing the whole merge field.\nconst suffixes = [\"t\", \"a\", \"r\"];\nfor (const suffix of suffixes) {\n  const placeholder = \"{{client_visits_\" + suffix + \"}}\";\n  const found = context.document.body.search(placeholder, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(placeholder, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: header cell \"Result up/Down\" -> \"Result\"\n$find = $d.Content.Find\n$find.Text = \"Result up/Down\"\n$find.Replacement.Text = \"Result\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# Change 2-4: collapse the split \"{{client_visits_X}}\" placeholders (each was\n# spread across three runs: \"{{\", \"client_visits_\", \"X}}\") back into a single\n# run containing the whole merge field. Searching/replacing with the same\n# literal text merges the matched range into one run.\n$suffixes = @(\"t\", \"a\", \"r\")\nforeach ($suffix in $suffixes) {\n    $placeholder = \"{{client_visits_$suffix}}\"\n    $find = $d.Content.Find\n    $find.Text = $placeholder\n    $find.Replacement.Text = $placeholder\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
